$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Renumber the client requirement IDs in column A (rows 23-25) so the
# duplicate "CRS_Client_006" at A23 becomes CRS_Client_007, and the
# following rows shift up accordingly, adding a new CRS_Client_009.
$ws.Range("A23").Value = "CRS_Client_007"
$ws.Range("A24").Value = "CRS_Client_008"
$ws.Range("A25").Value = "CRS_Client_009"

# Update the view: scroll/selection moved to B26 with no frozen top-left cell override.
$ws.Range("B26").Select()
